$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the entire first row (the header row) and delete it,
# shifting all the data rows up by one.
$ws.Rows.Item(1).Select()
$ws.Rows.Item(1).Delete()

# After deleting the row, Excel leaves the entire new row 1 selected.
$ws.Rows.Item(1).Select()
